# Add three new action items from the Oct 18 meeting to the bottom of the
# Action_Items sheet (rows 18-20), matching the style/format of the existing
# "white" rows (e.g. rows 14/15) and then move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 18 (Item 17): Upgrade Training tier to caArray 2.4.1. ---
$ws.Range("A14:E14").Copy()
$ws.Range("A18:E18").PasteSpecial(-4122)
$ws.Rows.Item(18).RowHeight = $ws.Rows.Item(14).RowHeight

$ws.Cells.Item(18, 1).Value = 17
$ws.Cells.Item(18, 2).Value = "Upgrade Training tier to caArray 2.4.1."
$ws.Cells.Item(18, 3).Value = "Don Swan"
$ws.Cells.Item(18, 4).Value = 39372
$ws.Cells.Item(18, 5).Value = "Not Started"

# --- Row 19 (Item 18): Upgrade Curation tier to caArray 2.4.1. ---
$ws.Range("A14:E14").Copy()
$ws.Range("A19:E19").PasteSpecial(-4122)
$ws.Rows.Item(19).RowHeight = $ws.Rows.Item(14).RowHeight

$ws.Cells.Item(19, 1).Value = 18
$ws.Cells.Item(19, 2).Value = "Upgrade Curation tier to caArray 2.4.1."
$ws.Cells.Item(19, 3).Value = "Quy Phung"
$ws.Cells.Item(19, 4).Value = 39372
$ws.Cells.Item(19, 5).Value = "Not Started"

# --- Row 20 (Item 19): Create wireframes re: permissions across caIntegrator/caArray ---
$ws.Range("A15:E15").Copy()
$ws.Range("A20:E20").PasteSpecial(-4122)
$ws.Rows.Item(20).RowHeight = $ws.Rows.Item(15).RowHeight

$ws.Cells.Item(20, 1).Value = 19
$ws.Cells.Item(20, 2).Value = "Create wireframes to depict how permissions will work across caIntegrator and caArray."
$ws.Cells.Item(20, 3).Value = "Andy Evans and Will Fitzhugh"
$ws.Cells.Item(20, 4).Value = 39372
$ws.Cells.Item(20, 5).Value = "Not Started"

# Clear the clipboard marquee and move the selection like the authored edit.
$excel.CutCopyMode = $false
$ws.PageSetup.Orientation = 1
$ws.Range("C19").Select() | Out-Null
